# "more training data and reduction of the features"
#
# For the Features sheet (Feuil1), three feature columns are removed:
#   H -> produces_colorless
#   N -> power
#   O -> Toughness
# Removing them shifts every column to its right one step left, which is
# exactly what Excel's "delete entire column" does (and also keeps the
# remaining columns' custom widths attached to their new positions).
#
# Delete from right to left so earlier column indices stay valid while we
# work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Column O (15) = "Toughness"
$ws.Columns.Item(15).EntireColumn.Delete()

# Column N (14) = "power"
$ws.Columns.Item(14).EntireColumn.Delete()

# Column H (8) = "produces_colorless"
$ws.Columns.Item(8).EntireColumn.Delete()

# The former "produces_colorless" column has now been replaced (shifted
# left) by what used to be column I ("Need_W"); select it the way the
# author apparently left the selection afterwards (whole column H).
[void]$ws.Columns.Item(8).Select()
